$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 313.875
$ws.Range("I2").Value = 315.85715
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 315.85715
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -202.85715
$ws.Range("N2").Value = -526
$ws.Range("H4").Value = 1891.8334
$ws.Range("I4").Value = 2070.4
$ws.Range("J4").Value = 999
$ws.Range("K4").Value = 2070.4
$ws.Range("L4").Value = 999
$ws.Range("M4").Value = -1956.4
$ws.Range("N4").Value = -1227
$ws.Range("H43").Value = 14434
$ws.Range("I43").Value = 14397.5
$ws.Range("K43").Value = 14397.5
$ws.Range("M43").Value = -14328.5
$ws.Range("H98").Value = 26064.654
$ws.Range("I98").Value = 43968.383
$ws.Range("J98").Value = 11517.875
$ws.Range("K98").Value = 43968.383
$ws.Range("L98").Value = 11517.875
$ws.Range("M98").Value = -42470.383
$ws.Range("N98").Value = -14513.875
$ws.Range("H113").Value = 13780.1
$ws.Range("I113").Value = 14884.5
$ws.Range("J113").Value = 12123.5
$ws.Range("K113").Value = 14884.5
$ws.Range("L113").Value = 12123.5
$ws.Range("M113").Value = -11630.5
$ws.Range("N113").Value = -18631.5
$ws.Range("H122").Value = 26064.654
$ws.Range("I122").Value = 43968.383
$ws.Range("J122").Value = 11517.875
$ws.Range("K122").Value = 131905.149
$ws.Range("L122").Value = 34553.625
$ws.Range("M122").Value = -129455.149
$ws.Range("N122").Value = -39453.625
$ws.Range("H135").Value = 6151.815
$ws.Range("I135").Value = 7155.095
$ws.Range("J135").Value = 2640.3333
$ws.Range("K135").Value = 64395.855
$ws.Range("L135").Value = 23762.9997
$ws.Range("M135").Value = -61860.855
$ws.Range("N135").Value = -28832.9997
$ws.Range("H137").Value = 10903.458
$ws.Range("I137").Value = 15583.866
$ws.Range("K137").Value = 46751.598
$ws.Range("M137").Value = -44201.598
$ws.Range("H138").Value = 282305.72
$ws.Range("I138").Value = 488940.56
$ws.Range("J138").Value = 3348.7
$ws.Range("K138").Value = 1466821.68
$ws.Range("L138").Value = 10046.1
$ws.Range("M138").Value = -1461681.68
$ws.Range("N138").Value = -20326.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3731.9565
$ws.Range("I2").Value = 3468.611
$ws.Range("K2").Value = 3468.611
$ws.Range("M2").Value = -3355.611
$ws.Range("H32").Value = 4060.7407
$ws.Range("I32").Value = 3793.7793
$ws.Range("J32").Value = 9199.75
$ws.Range("K32").Value = 3793.7793
$ws.Range("L32").Value = 9199.75
$ws.Range("M32").Value = -3506.7793
$ws.Range("N32").Value = -9773.75
$ws.Range("H45").Value = 295283
$ws.Range("I45").Value = 677293.7
$ws.Range("J45").Value = 8775
$ws.Range("K45").Value = 677293.7
$ws.Range("L45").Value = 8775
$ws.Range("M45").Value = -676916.7
$ws.Range("N45").Value = -9529
$ws.Range("H97").Value = 6064433
$ws.Range("I97").Value = 4910.0835
$ws.Range("J97").Value = 22223160
$ws.Range("K97").Value = 4910.0835
$ws.Range("L97").Value = 22223160
$ws.Range("M97").Value = -4414.0835
$ws.Range("N97").Value = -22224152
$ws.Range("H116").Value = 3731.9565
$ws.Range("I116").Value = 3468.611
$ws.Range("K116").Value = 3468.611
$ws.Range("M116").Value = -1174.611
$ws.Range("H122").Value = 1207507.2
$ws.Range("I122").Value = 7263.2856
$ws.Range("K122").Value = 21789.8568
$ws.Range("M122").Value = -19339.8568
$ws.Range("H132").Value = 2285.283
$ws.Range("I132").Value = 1315.8572
$ws.Range("J132").Value = 4170.278
$ws.Range("K132").Value = 3947.5716
$ws.Range("L132").Value = 12510.834
$ws.Range("M132").Value = -1417.5716
$ws.Range("N132").Value = -17570.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3731.9565
$ws.Range("I3").Value = 3468.611
$ws.Range("K3").Value = 3468.611
$ws.Range("M3").Value = -3354.611
$ws.Range("H20").Value = 4936.8
$ws.Range("I20").Value = 2820.6667
$ws.Range("J20").Value = 5843.7144
$ws.Range("K20").Value = 2820.6667
$ws.Range("L20").Value = 5843.7144
$ws.Range("M20").Value = -2573.6667
$ws.Range("N20").Value = -6337.7144
$ws.Range("H134").Value = 6717.2905
$ws.Range("I134").Value = 8214.75
$ws.Range("J134").Value = 3994.6365
$ws.Range("K134").Value = 24644.25
$ws.Range("L134").Value = 11983.9095
$ws.Range("M134").Value = -22109.25
$ws.Range("N134").Value = -17053.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 170011.03
$ws.Range("I99").Value = 314755.2
$ws.Range("J99").Value = 4589.143
$ws.Range("K99").Value = 314755.2
$ws.Range("L99").Value = 4589.143
$ws.Range("M99").Value = -313257.2
$ws.Range("N99").Value = -7585.143
$ws.Range("H122").Value = 8963.200000000001
$ws.Range("I122").Value = 8963.200000000001
$ws.Range("K122").Value = 26889.6
$ws.Range("M122").Value = -24439.6
$ws.Range("H126").Value = 170011.03
$ws.Range("I126").Value = 314755.2
$ws.Range("J126").Value = 4589.143
$ws.Range("K126").Value = 944265.6000000001
$ws.Range("L126").Value = 13767.429
$ws.Range("M126").Value = -941795.6000000001
$ws.Range("N126").Value = -18707.429
$ws.Range("H132").Value = 1838.2222
$ws.Range("I132").Value = 1873
$ws.Range("J132").Value = 1716.5
$ws.Range("K132").Value = 5619
$ws.Range("L132").Value = 5149.5
$ws.Range("M132").Value = -3089
$ws.Range("N132").Value = -10209.5
$ws.Range("H134").Value = 5686.522
$ws.Range("I134").Value = 6247.55
$ws.Range("J134").Value = 1946.3334
$ws.Range("K134").Value = 18742.65
$ws.Range("L134").Value = 5839.0002
$ws.Range("M134").Value = -16207.65
$ws.Range("N134").Value = -10909.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1614.6666
$ws.Range("I86").Value = 444
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 1332
$ws.Range("L86").Value = 6600
$ws.Range("M86").Value = -146
$ws.Range("N86").Value = -8972
$ws.Range("H89").Value = 1614.6666
$ws.Range("I89").Value = 444
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 3996
$ws.Range("L89").Value = 19800
$ws.Range("M89").Value = 1932
$ws.Range("N89").Value = -31656
$ws.Range("H121").Value = 2208
$ws.Range("I121").Value = 525
$ws.Range("K121").Value = 1575
$ws.Range("M121").Value = -265
$ws.Range("H131").Value = 1961.4597
$ws.Range("J131").Value = 2028.4684
$ws.Range("L131").Value = 6085.4052
$ws.Range("N131").Value = -16165.4052
$ws.Range("H140").Value = 3071.7058
$ws.Range("I140").Value = 2925.7334
$ws.Range("J140").Value = 4166.5
$ws.Range("K140").Value = 8777.200199999999
$ws.Range("L140").Value = 12499.5
$ws.Range("M140").Value = -3597.200199999999
$ws.Range("N140").Value = -22859.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15937.909
$ws.Range("I102").Value = 16931.7
$ws.Range("K102").Value = 16931.7
$ws.Range("M102").Value = -15309.7
$ws.Range("H122").Value = 12945.667
$ws.Range("I122").Value = 8240.538
$ws.Range("K122").Value = 24721.614
$ws.Range("M122").Value = -22271.614
$ws.Range("H132").Value = 4837.4375
$ws.Range("I132").Value = 4999.852
$ws.Range("J132").Value = 3960.4
$ws.Range("K132").Value = 14999.556
$ws.Range("L132").Value = 11881.2
$ws.Range("M132").Value = -12469.556
$ws.Range("N132").Value = -16941.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2802.6924
$ws.Range("I82").Value = 4516.3335
$ws.Range("J82").Value = 1333.8572
$ws.Range("K82").Value = 4516.3335
$ws.Range("L82").Value = 1333.8572
$ws.Range("M82").Value = -4155.3335
$ws.Range("N82").Value = -2055.8572
$ws.Range("H85").Value = 2802.6924
$ws.Range("I85").Value = 4516.3335
$ws.Range("J85").Value = 1333.8572
$ws.Range("K85").Value = 4516.3335
$ws.Range("L85").Value = 1333.8572
$ws.Range("M85").Value = -3268.3335
$ws.Range("N85").Value = -3829.8572
$ws.Range("H122").Value = 6388.1924
$ws.Range("I122").Value = 6068.9375
$ws.Range("K122").Value = 18206.8125
$ws.Range("M122").Value = -15756.8125
$ws.Range("H132").Value = 395890.4
$ws.Range("I132").Value = 623960.4399999999
$ws.Range("J132").Value = 4913.143
$ws.Range("K132").Value = 1871881.32
$ws.Range("L132").Value = 14739.429
$ws.Range("M132").Value = -1869351.32
$ws.Range("N132").Value = -19799.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4378.1836
$ws.Range("I122").Value = 2826.3794
$ws.Range("J122").Value = 6628.3
$ws.Range("K122").Value = 8479.138199999999
$ws.Range("L122").Value = 19884.9
$ws.Range("M122").Value = -6029.138199999999
$ws.Range("N122").Value = -24784.9
$ws.Range("H126").Value = 15046.343
$ws.Range("I126").Value = 18995.809
$ws.Range("K126").Value = 56987.427
$ws.Range("M126").Value = -54517.427
$ws.Range("H132").Value = 3327.625
$ws.Range("I132").Value = 2604.9744
$ws.Range("J132").Value = 4985.4707
$ws.Range("K132").Value = 7814.9232
$ws.Range("L132").Value = 14956.4121
$ws.Range("M132").Value = -5284.9232
$ws.Range("N132").Value = -20016.4121
